$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new rows (list is kept sorted alphabetically by column A) ---
# New row 8: "Divorce" (sorts between "Criminal Court fee waiver" and "E-filing exemption")
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "Divorce"
$ws.Range("B8").Value = "https://www.illinoislegalaid.org/legal-information/divorce"

# New row 26: "Parenting plan" (sorts between "Order of protection" and "Personnel file request")
$ws.Rows.Item(26).Insert()
$ws.Range("A26").Value = "Parenting plan"
$ws.Range("B26").Value = "https://www.illinoislegalaid.org/legal-information/parenting-plan"

# --- Rebuild hyperlinks so each one anchors to its (now shifted) row ---
# Inserting rows does not relocate pre-existing hyperlink anchors automatically,
# so clear every hyperlink on the sheet and re-add them against the new layout.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.illinoislegalaid.org/legal-information/appearance")
$ws.Hyperlinks.Add($ws.Range("B16"), "https://www.illinoislegalaid.org/legal-information/fee-waiver")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://www.illinoislegalaid.org/legal-information/collection-proof-debtor-letter")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.illinoislegalaid.org/legal-information/request-collection-agency-stop-contacting")
$ws.Hyperlinks.Add($ws.Range("B13"), "https://www.illinoislegalaid.org/legal-information/end-illegal-lockout-demand")
$ws.Hyperlinks.Add($ws.Range("B40"), "https://www.illinoislegalaid.org/legal-information/security-deposit-demand-letter")
$ws.Hyperlinks.Add($ws.Range("B18"), "https://www.illinoislegalaid.org/legal-information/housing-discrimination-complaint-idhr")
$ws.Hyperlinks.Add($ws.Range("B44"), "https://www.illinoislegalaid.org/legal-information/stop-wage-assignment-letter")
$ws.Hyperlinks.Add($ws.Range("B34"), "https://www.illinoislegalaid.org/legal-information/request-time-work-due-domestic-abuse-letter")
$ws.Hyperlinks.Add($ws.Range("B9"), "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-circuit-court")
$ws.Hyperlinks.Add($ws.Range("B10"), "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-appellate-court")
$ws.Hyperlinks.Add($ws.Range("B11"), "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-supreme-court")
$ws.Hyperlinks.Add($ws.Range("B36"), "https://www.illinoislegalaid.org/legal-information/respond-lawsuit")
$ws.Hyperlinks.Add($ws.Range("B47"), "https://www.illinoislegalaid.org/legal-information/voluntary-acknowledgment-parentage-vap")
$ws.Hyperlinks.Add($ws.Range("B19"), "https://www.illinoislegalaid.org/legal-information/interpreter-request")
$ws.Hyperlinks.Add($ws.Range("B21"), "https://www.illinoislegalaid.org/legal-information/motion")
$ws.Hyperlinks.Add($ws.Range("B45"), "https://www.illinoislegalaid.org/legal-information/transfer-death-instrument-or-todi")
$ws.Hyperlinks.Add($ws.Range("B30"), "https://www.illinoislegalaid.org/legal-information/power-attorney-agent-resign-letter")
$ws.Hyperlinks.Add($ws.Range("B31"), "https://www.illinoislegalaid.org/legal-information/power-attorney-revocation")
$ws.Hyperlinks.Add($ws.Range("B29"), "https://www.illinoislegalaid.org/legal-information/power-attorney-property")
$ws.Hyperlinks.Add($ws.Range("B28"), "https://www.illinoislegalaid.org/legal-information/power-attorney-health-care")
$ws.Hyperlinks.Add($ws.Range("B25"), "https://www.illinoislegalaid.org/legal-information/order-protection")
$ws.Hyperlinks.Add($ws.Range("B23"), "https://www.illinoislegalaid.org/legal-information/name-change-adult")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://www.illinoislegalaid.org/legal-information/cannabis-expungement")
$ws.Hyperlinks.Add($ws.Range("B12"), "https://www.illinoislegalaid.org/legal-information/emergency-order-protection-cook-county")
$ws.Hyperlinks.Add($ws.Range("B41"), "https://www.illinoislegalaid.org/legal-information/short-term-guardian-appointment")
$ws.Hyperlinks.Add($ws.Range("B33"), "https://www.illinoislegalaid.org/legal-information/remove-eviction-public-record")
$ws.Hyperlinks.Add($ws.Range("B38"), "https://www.illinoislegalaid.org/legal-information/respond-eviction")
$ws.Hyperlinks.Add($ws.Range("B42"), "https://www.illinoislegalaid.org/legal-information/small-claims-complaint")
$ws.Hyperlinks.Add($ws.Range("B39"), "https://www.illinoislegalaid.org/legal-information/security-deposit-complaint")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://www.illinoislegalaid.org/legal-information/criminal-court-fee-waiver")
$ws.Hyperlinks.Add($ws.Range("B17"), "https://www.illinoislegalaid.org/legal-information/financial-affidavit")
$ws.Hyperlinks.Add($ws.Range("B35"), "https://www.illinoislegalaid.org/legal-information/special-process-server-request")
$ws.Hyperlinks.Add($ws.Range("B46"), "https://www.illinoislegalaid.org/legal-information/vacate-default-judgment-within-30-days")
$ws.Hyperlinks.Add($ws.Range("B22"), "https://www.illinoislegalaid.org/legal-information/motion-continue-or-extend-time")

# --- Re-apply the plain "Hyperlink" cell style to every row in column B ---
# (Hyperlinks.Add() registers a second, duplicate "Hyperlink" style the first time
# it runs; resetting .Style here keeps every data row on the original style index.)
for ($r = 2; $r -le 47; $r++) {
    $ws.Range("B" + $r).Style = "Hyperlink"
}

# --- Restore the selection shown in the sheet view ---
$ws.Range("A27").Select() | Out-Null

